# Apply the 16.9.1 indicator workbook edits.
# The source commit removes the separate Russian/Kyrgyz "Sex"/"Male"/"Female"/
# "Areas"/"Age (in month)"/"Education" style header labels and replaces them
# with new, more specific header texts (e.g. "Sex" -> "By sex",
# "Male" -> "Men", "Areas" -> "By territory", etc.) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matches the order new shared strings were introduced in the
# canonical edit, so the rebuilt shared-strings table lines up cell-for-cell.
$ws.Range("C6").Value  = "By sex"
$ws.Range("C12").Value = "By territory"
$ws.Range("C22").Value = "By age (in month)"
$ws.Range("C28").Value = "Education of mother"
$ws.Range("C29").Value = "Preschool or not /primary"
$ws.Range("C30").Value = "Basic general"
$ws.Range("C31").Value = "Average total"
$ws.Range("C32").Value = "Vocational primary /secondary"
$ws.Range("C33").Value = "Higher"
$ws.Range("C7").Value  = "Men"
$ws.Range("C8").Value  = "Woman"
$ws.Range("B7").Value  = "Мужчины"
$ws.Range("B8").Value  = "Женщины"
$ws.Range("A7").Value  = "Эркектер"
$ws.Range("A8").Value  = "Аялдар"
$ws.Range("A22").Value = "Жаш курагы боюнча (айларда)"
$ws.Range("B22").Value = "По возрасту (в месяцах)"
